# 11.8.18 plots and analysis
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row fixes ---
$ws.Range("E3").Value = "Report transformed coefficient? (OR/RR)"

# --- Fix typo in row 7 note ---
$ws.Range("M7").Value = "Figure shows predicted probs, but not described in text"

# --- Row 9: fill in missing J9 value ---
$ws.Range("J9").Value = "N"

# --- Row 10: fill in full result row ---
$ws.Range("A10").Value = "JAbP"
$ws.Range("C10").Value = "Y"
$ws.Range("D10").Value = "LR"
$ws.Range("E10").Value = "Y (OR)"
$ws.Range("F10").Value = "N"
$ws.Range("G10").Value = "N"
$ws.Range("H10").Value = "Y"
$ws.Range("I10").Value = "N"
$ws.Range("J10").Value = "N"
$ws.Range("K10").Value = "N"
$ws.Range("L10").Value = 195

# --- Row 11: fill in full result row ---
$ws.Range("A11").Value = "JCCP"
$ws.Range("C11").Value = "Y"
$ws.Range("D11").Value = "LR"
$ws.Range("E11").Value = "Y (OR)"
$ws.Range("F11").Value = "N"
$ws.Range("G11").Value = "N"
$ws.Range("H11").Value = "N"
$ws.Range("I11").Value = "N"
$ws.Range("J11").Value = "N"
$ws.Range("K11").Value = "N"
$ws.Range("L11").Value = 307

# --- Row 12: methods-paper style row ---
$ws.Range("C12").Value = "N"
$ws.Range("D12").Value = "-"
$ws.Range("E12").Value = "-"
$ws.Range("F12").Value = "-"
$ws.Range("G12").Value = "-"
$ws.Range("H12").Value = "-"
$ws.Range("I12").Value = "-"
$ws.Range("J12").Value = "-"
$ws.Range("K12").Value = "-"
$ws.Range("L12").Value = "-"
$ws.Range("M12").Value = "This is David Atkins' tutorial paper"

# --- Update current selection to match saved view ---
$ws.Range("A12").Select()
